# Updated cryptos list - apply price/volume changes scraped on refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.233.70'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '1.702.27'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.17'
$ws.Range("E5").Value = '  -1.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5304'
$ws.Range("E6").Value = '  -1.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2657'
$ws.Range("E8").Value = '  -0.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06581'
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.73'
$ws.Range("E10").Value = '  -4.88%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07624'
$ws.Range("E11").Value = '  -1.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.493'
$ws.Range("E12").Value = '  -3.13%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.705.21'
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.938.70'
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5791'
$ws.Range("E15").Value = '  -1.43%  '
$ws.Range("D16").Value = '0.0₅8147'
$ws.Range("E16").Value = '  -1.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.36'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").Value = '27.278.58'
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.34'
$ws.Range("E19").Value = '  -3.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.609'
$ws.Range("E21").Value = '  -2.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.32'
$ws.Range("E22").Value = '  -3.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.967'
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.95'
$ws.Range("E25").Value = '  -2.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.701'
$ws.Range("E26").Value = '  +0.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1195'
$ws.Range("E27").Value = '  -3.11%  '
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '16.14'
$ws.Range("E29").Value = '  -3.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05365'
$ws.Range("E30").Value = '  -3.39%  '
$ws.Range("E31").Value = '  -1.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.461'
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.395'
$ws.Range("E33").Value = '  -2.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.638'
$ws.Range("E34").Value = '  -1.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.861'
$ws.Range("E35").Value = '  +1.55%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9436'
$ws.Range("E36").Value = '  -1.80%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.405'
$ws.Range("E37").Value = '  -1.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5825'
$ws.Range("E38").Value = '  -2.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01632'
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.779'
$ws.Range("E40").Value = '  -1.63%  '
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '1.039.25'
$ws.Range("E42").Value = '  -1.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8396'
$ws.Range("E43").Value = '  -2.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.70'
$ws.Range("E44").Value = '  -0.95%  '
$ws.Range("D45").Value = '1.845.88'
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.66'
$ws.Range("E47").Value = '  -2.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4516'
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.008'
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.029'
$ws.Range("E50").Value = '  -2.46%  '
$ws.Range("E51").Value = '  -0.86%  '
